# Apply the "feat: add 2022-Q3 data" change:
# 1. Insert a new row into "总计" (summary) sheet for 2022-Q3.
# 2. Insert a new worksheet "2022-Q3" (before "2022-Q2") with the quarter's fund data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" summary sheet: insert a new row right after the header
#    and fill it in with the 2022-Q3 totals, pushing existing rows down.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Re-apply the same formatting used by the other data rows (copied from what
# is now row 3, the former row 2) so the new row matches the sheet's style.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 0.61

# ---------------------------------------------------------------------------
# 2) Create the new "2022-Q3" worksheet (positioned right before "2022-Q2"),
#    reusing "2022-Q2"'s layout/formatting, then fill in this quarter's data.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Make sure the fund-code / numeric-looking text columns stay text so values
# like "002367" or "8.50" are not silently turned into numbers.
$q3.Range("B2:B7").NumberFormat = "@"
$q3.Range("D2:G7").NumberFormat = "@"

# Row 2 - 257010 国联安小盘精选混合 (values updated)
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "257010"
$q3.Range("C2").Value = "国联安小盘精选混合"
$q3.Range("D2").Value = "8.50"
$q3.Range("E2").Value = "74.70"
$q3.Range("F2").Value = "5.88"
$q3.Range("G2").Value = "0.4998"
$q3.Range("H2").Value = 2

# Row 3 - 159617 华夏中证智选500价值稳健策略ETF (new fund in this slot)
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "159617"
$q3.Range("C3").Value = "华夏中证智选500价值稳健策略ETF"
$q3.Range("D3").Value = "2.93"
$q3.Range("E3").Value = "97.05"
$q3.Range("F3").Value = "1.50"
$q3.Range("G3").Value = "0.0440"
$q3.Range("H3").Value = 2

# Row 4 - 006138 国联安价值优选股票 (values updated)
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "006138"
$q3.Range("C4").Value = "国联安价值优选股票"
$q3.Range("D4").Value = "0.57"
$q3.Range("E4").Value = "94.64"
$q3.Range("F4").Value = "5.83"
$q3.Range("G4").Value = "0.0332"
$q3.Range("H4").Value = 6

# Rows 5-7 are brand-new this quarter; copy formatting from row 4 first.
$q3.Range("A4:H4").Copy()
$q3.Range("A5:H5").PasteSpecial(-4122) # xlPasteFormats
$q3.Range("A6:H6").PasteSpecial(-4122)
$q3.Range("A7:H7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q3.Range("B5:B7").NumberFormat = "@"
$q3.Range("D5:G7").NumberFormat = "@"

# Row 5 - 002367 国联安安稳灵活配置混合
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "002367"
$q3.Range("C5").Value = "国联安安稳灵活配置混合"
$q3.Range("D5").Value = "0.57"
$q3.Range("E5").Value = "47.79"
$q3.Range("F5").Value = "3.46"
$q3.Range("G5").Value = "0.0197"
$q3.Range("H5").Value = 3

# Row 6 - 004351 汇丰晋信珠三角区域发展混合
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "004351"
$q3.Range("C6").Value = "汇丰晋信珠三角区域发展混合"
$q3.Range("D6").Value = "0.42"
$q3.Range("E6").Value = "93.94"
$q3.Range("F6").Value = "2.81"
$q3.Range("G6").Value = "0.0118"
$q3.Range("H6").Value = 9

# Row 7 - 004250 银河量化优选混合
$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "004250"
$q3.Range("C7").Value = "银河量化优选混合"
$q3.Range("D7").Value = "0.24"
$q3.Range("E7").Value = "86.30"
$q3.Range("F7").Value = "1.72"
$q3.Range("G7").Value = "0.0041"
$q3.Range("H7").Value = 8

# Restore the originally-active tab (the last sheet, 2020-Q4), since adding
# and editing sheets above moved Excel's "active sheet" selection.
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Host "Done applying 2022-Q3 update"
